# Andrew Hickman Resume - content edits
# 1) Rewrite the "Summary" paragraph.
# 2) Re-point the "_GoBack" bookmark to the end of the (new) Summary text
#    (Word automatically relocates _GoBack to the site of the most recent
#    edit; _gjdgxs keeps its position but is renumbered as a side effect).
# 3) Reword the Skyline RV Park Construction date range.

$d = $word.ActiveDocument

# --- 1) Summary paragraph -------------------------------------------------
$oldSummary = "My interests in hardware, programming, and network design led me on" + `
  " the path to being an Electronic Systems Engineering major. Through ESET," + `
  " I was introduced to programming and have strived to further my knowledge" + `
  " on the subject so that I can be a valued Software Engineering and Project" + `
  " Manager for my Senior Capstone Design team."

$newSummary = "My strong leadership skills shine bright in stressful situations. " + `
  "Working at Internet2 Technology Evaluation Center has taught me the values " + `
  "in being able to independently troubleshoot issues to improve my problem " + `
  "solving skills. I work hard to improve my programming expertise so that I " + `
  "can be a valued Software Engineer and Project Manager for my Senior " + `
  "Capstone Design team."

$summaryPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text.StartsWith("My interests in hardware")) {
        $summaryPara = $p
        break
    }
}

$summaryRange = $summaryPara.Range
$found = $summaryRange.Find.Execute($oldSummary, $true, $false, $false, $false, $false, `
  $true, 1, $false, $newSummary, 2)

# --- 2) Move the _GoBack bookmark to the end of the Summary paragraph -----
$summaryPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text.StartsWith("My strong leadership skills")) {
        $summaryPara = $p
        break
    }
}
$endOfSummaryText = $summaryPara.Range.End - 1
$d.Bookmarks.Add("_GoBack", $d.Range($endOfSummaryText, $endOfSummaryText))

# --- 3) Skyline RV Park Construction date range ----------------------------
$skylinePara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "*Skyline RV Park Construction*") {
        $skylinePara = $p
        break
    }
}

$endash = [char]8211
$square = [char]9642
$oldDate = "(June 2017 $endash August 2017) $square "
$newDate = "(Summer 2016 & 2017) $square "
$skylinePara.Range.Find.Execute($oldDate, $true, $false, $false, $false, $false, `
  $true, 1, $false, $newDate, 2) | Out-Null
